# Staging.Project.xlsx — "moved staging files StagingTemplates directory"
#
# The underlying change (per the OOXML diff) renames the four "...SourceKey"
# staging columns to "...BusinessKey" (reflecting the move/rename of the
# staging source tables) and leaves every other header text as-is. The
# shared-string table in the target file also happens to come out sorted
# alphabetically (a side effect of however the source tool regenerated the
# sheet) and a couple of purely cosmetic, non-data attributes changed
# (worksheet codeName, the saved window size, and the explicit column-width
# overrides were dropped back to sheet defaults).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths: the explicit bestFit/customWidth overrides on columns
# B:K are dropped in the target (only column A keeps its custom width), i.e.
# columns B:K revert to the sheet's default width. Deleting & shifting the
# columns back out clears their stored width overrides; do this FIRST, while
# the header row still holds its original text, so the delete's left-shift
# doesn't disturb anything we still need to read.
$ws.Range("B1:K1").EntireColumn.Delete()

# --- Header row (row 2): rename the four SourceKey columns to BusinessKey.
# The column delete above also removed the row-2 header cells along with the
# column formatting, so re-enter the full header set (text unchanged except
# for the four renamed columns).
$ws.Range("A2").Value = "Code"
$ws.Range("B2").Value = "LongName"
$ws.Range("C2").Value = "OutcomeBusinessKey"
$ws.Range("D2").Value = "ProgrammeBusinessKey"
$ws.Range("E2").Value = "ProjectID"
$ws.Range("F2").Value = "ProjectParentID"
$ws.Range("G2").Value = "ProjectSiteName"
$ws.Range("H2").Value = "SectorBusinessKey"
$ws.Range("I2").Value = "ShortName"
$ws.Range("J2").Value = "SubSectorBusinessKey"
$ws.Range("K2").Value = "TextDescription"

# --- Worksheet code name: Sheet42 -> Sheet44.
$ws.CodeName = "Sheet44"

# --- Saved window size (bookViews/workbookView).
$wb.Windows.Item(1).Width = 28800
$wb.Windows.Item(1).Height = 12585
